# Commit: "Add experimental setup and fix diagrams"
#
# 1) Refresh the cached "datetimeFigureOut" footer/date field
#    (08/05/2017 -> 11/05/2017) on the slide master and on every slide
#    layout.
# 2) Un-bold the four small connector labels in the diagram group on
#    slide 1 ("utilities", "dpdgraph++", "file-graph", "library of
#    queries").

$p = $ppt.ActivePresentation

$newDate = "11/05/2017"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch { }
        if ($isDate -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 1) Date placeholder on the slide master -------------------------------
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- Date placeholder on every slide layout ---------------------------------
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# --- 2) Remove bold from the four diagram labels ----------------------------
$labels = @("utilities", "dpdgraph++", "file-graph", "library of queries")

$slide = $p.Slides.Item(1)
for ($si = 1; $si -le $slide.Shapes.Count; $si++) {
    $top = $slide.Shapes.Item($si)
    if ($top.Type -eq 6) {
        # msoGroup - walk the group items looking for our labels
        for ($gi = 1; $gi -le $top.GroupItems.Count; $gi++) {
            $item = $top.GroupItems.Item($gi)
            if ($item.HasTextFrame) {
                $tr = $item.TextFrame.TextRange
                if ($labels -contains $tr.Text) {
                    $tr.Font.Bold = $false
                }
            }
        }
    } elseif ($top.HasTextFrame) {
        $tr = $top.TextFrame.TextRange
        if ($labels -contains $tr.Text) {
            $tr.Font.Bold = $false
        }
    }
}
